$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-13 with new model metrics (reordered + new values)
# Row 2
$ws.Cells.Item(2, 1).Value = "model_12_5_2"
$ws.Cells.Item(2, 2).Value = 0.9452851570346646
$ws.Cells.Item(2, 3).Value = 0.6874607407731951
$ws.Cells.Item(2, 4).Value = 0.7562110000974103
$ws.Cells.Item(2, 5).Value = 0.9562878173996451
$ws.Cells.Item(2, 6).Value = 0.921363701105245
$ws.Cells.Item(2, 7).Value = 0.3658784068478836
$ws.Cells.Item(2, 8).Value = 2.089951465560607
$ws.Cells.Item(2, 9).Value = 0.7165842306809362
$ws.Cells.Item(2, 10).Value = 0.546092905599731
$ws.Cells.Item(2, 11).Value = 0.6313385681403335
$ws.Cells.Item(2, 12).Value = 0.8889694374839107
$ws.Cells.Item(2, 13).Value = 0.6048788365018928
$ws.Cells.Item(2, 14).Value = 1.029181249581512
$ws.Cells.Item(2, 15).Value = 0.6306297808435866
$ws.Cells.Item(2, 16).Value = 140.0109084450115
$ws.Cells.Item(2, 17).Value = 224.1133403609173

# Row 3
$ws.Cells.Item(3, 1).Value = "model_12_5_3"
$ws.Cells.Item(3, 2).Value = 0.9458206739679796
$ws.Cells.Item(3, 3).Value = 0.6870406628782313
$ws.Cells.Item(3, 4).Value = 0.6968655993997701
$ws.Cells.Item(3, 5).Value = 0.9562014228675083
$ws.Cells.Item(3, 6).Value = 0.9104329541171079
$ws.Cells.Item(3, 7).Value = 0.3622974026489774
$ws.Cells.Item(3, 8).Value = 2.09276052837851
$ws.Cells.Item(3, 9).Value = 0.8910218727417447
$ws.Cells.Item(3, 10).Value = 0.5471722257863676
$ws.Cells.Item(3, 11).Value = 0.719097049264056
$ws.Cells.Item(3, 12).Value = 0.8427768239223999
$ws.Cells.Item(3, 13).Value = 0.6019114574827242
$ws.Cells.Item(3, 14).Value = 1.028895640550411
$ws.Cells.Item(3, 15).Value = 0.6275360743562506
$ws.Cells.Item(3, 16).Value = 140.0305797003368
$ws.Cells.Item(3, 17).Value = 224.1330116162427

# Row 4
$ws.Cells.Item(4, 1).Value = "model_12_5_1"
$ws.Cells.Item(4, 2).Value = 0.9436083522572771
$ws.Cells.Item(4, 3).Value = 0.6867359121038118
$ws.Cells.Item(4, 4).Value = 0.8132562913741591
$ws.Cells.Item(4, 5).Value = 0.956931540807362
$ws.Cells.Item(4, 6).Value = 0.9323070163750328
$ws.Cells.Item(4, 7).Value = 0.3770912081152499
$ws.Cells.Item(4, 8).Value = 2.094798398210302
$ws.Cells.Item(4, 9).Value = 0.5489074438699957
$ws.Cells.Item(4, 10).Value = 0.5380509190135971
$ws.Cells.Item(4, 11).Value = 0.5434791814417963
$ws.Cells.Item(4, 12).Value = 0.9393771070439899
$ws.Cells.Item(4, 13).Value = 0.6140775261440935
$ws.Cells.Item(4, 14).Value = 1.030075545462785
$ws.Cells.Item(4, 15).Value = 0.6402200777477683
$ws.Cells.Item(4, 16).Value = 139.9505363789092
$ws.Cells.Item(4, 17).Value = 224.0529682948151

# Row 5
$ws.Cells.Item(5, 1).Value = "model_12_5_4"
$ws.Cells.Item(5, 2).Value = 0.9454404129055303
$ws.Cells.Item(5, 3).Value = 0.6856664978037459
$ws.Cells.Item(5, 4).Value = 0.6357323523052754
$ws.Cells.Item(5, 5).Value = 0.956552754783059
$ws.Cells.Item(5, 6).Value = 0.8995154942693158
$ws.Cells.Item(5, 7).Value = 0.3648402101245173
$ws.Cells.Item(5, 8).Value = 2.101949576559043
$ws.Cells.Item(5, 9).Value = 1.070714643357892
$ws.Cells.Item(5, 10).Value = 0.5427830634252211
$ws.Cells.Item(5, 11).Value = 0.8067488533915567
$ws.Cells.Item(5, 12).Value = 0.8002441203514609
$ws.Cells.Item(5, 13).Value = 0.6040200411613156
$ws.Cells.Item(5, 14).Value = 1.029098446450384
$ws.Cells.Item(5, 15).Value = 0.6297344248074095
$ws.Cells.Item(5, 16).Value = 140.0165916034686
$ws.Cells.Item(5, 17).Value = 224.1190235193744

# Row 6
$ws.Cells.Item(6, 1).Value = "model_12_5_0"
$ws.Cells.Item(6, 2).Value = 0.9404647133739481
$ws.Cells.Item(6, 3).Value = 0.684581253903157
$ws.Cells.Item(6, 4).Value = 0.8671468937153486
$ws.Cells.Item(6, 5).Value = 0.9581864156176484
$ws.Cells.Item(6, 6).Value = 0.9431483394383271
$ws.Cells.Item(6, 7).Value = 0.3981127358031974
$ws.Cells.Item(6, 8).Value = 2.109206607519369
$ws.Cells.Item(6, 9).Value = 0.3905034312401244
$ws.Cells.Item(6, 10).Value = 0.5223738653743756
$ws.Cells.Item(6, 11).Value = 0.45643864830725
$ws.Cells.Item(6, 12).Value = 0.9947392417087666
$ws.Cells.Item(6, 13).Value = 0.6309617546279628
$ws.Cells.Item(6, 14).Value = 1.031752152867228
$ws.Cells.Item(6, 15).Value = 0.6578231027933671
$ws.Cells.Item(6, 16).Value = 139.8420401160402
$ws.Cells.Item(6, 17).Value = 223.9444720319461

# Row 7
$ws.Cells.Item(7, 1).Value = "model_12_5_5"
$ws.Cells.Item(7, 2).Value = 0.9443012828011046
$ws.Cells.Item(7, 3).Value = 0.6834743833347401
$ws.Cells.Item(7, 4).Value = 0.5731984972081658
$ws.Cells.Item(7, 5).Value = 0.9571716077386564
$ws.Cells.Item(7, 6).Value = 0.8885497834204997
$ws.Cells.Item(7, 7).Value = 0.372457578378024
$ws.Cells.Item(7, 8).Value = 2.116608256107043
$ws.Cells.Item(7, 9).Value = 1.254524308536306
$ws.Cells.Item(7, 10).Value = 0.5350517814677196
$ws.Cells.Item(7, 11).Value = 0.8947880450020128
$ws.Cells.Item(7, 12).Value = 0.760942014413457
$ws.Cells.Item(7, 13).Value = 0.6102930266503329
$ws.Cells.Item(7, 14).Value = 1.029705982506077
$ws.Cells.Item(7, 15).Value = 0.6362744642755644
$ws.Cells.Item(7, 16).Value = 139.9752642623934
$ws.Cells.Item(7, 17).Value = 224.0776961782992

# Row 8
$ws.Cells.Item(8, 1).Value = "model_12_5_6"
$ws.Cells.Item(8, 2).Value = 0.9425140459853641
$ws.Cells.Item(8, 3).Value = 0.6805703115951535
$ws.Cells.Item(8, 4).Value = 0.5096533921080226
$ws.Cells.Item(8, 5).Value = 0.9578527306930218
$ws.Cells.Item(8, 6).Value = 0.8774474049788219
$ws.Cells.Item(8, 7).Value = 0.3844088391943493
$ws.Cells.Item(8, 8).Value = 2.13602779720168
$ws.Cells.Item(8, 9).Value = 1.441306403995576
$ws.Cells.Item(8, 10).Value = 0.5265425652471357
$ws.Cells.Item(8, 11).Value = 0.983924484621356
$ws.Cells.Item(8, 12).Value = 0.7245677630330192
$ws.Cells.Item(8, 13).Value = 0.6200071283415614
$ws.Cells.Item(8, 14).Value = 1.030659175474472
$ws.Cells.Item(8, 15).Value = 0.6464021153867511
$ws.Cells.Item(8, 16).Value = 139.912097214736
$ws.Cells.Item(8, 17).Value = 224.0145291306418

# Row 9
$ws.Cells.Item(9, 1).Value = "model_12_5_7"
$ws.Cells.Item(9, 2).Value = 0.9401611954998844
$ws.Cells.Item(9, 3).Value = 0.6770479611848401
$ws.Cells.Item(9, 4).Value = 0.4455530521337409
$ws.Cells.Item(9, 5).Value = 0.9583761289211308
$ws.Cells.Item(9, 6).Value = 0.866120673350424
$ws.Cells.Item(9, 7).Value = 0.4001423612246324
$ws.Cells.Item(9, 8).Value = 2.159581770614378
$ws.Cells.Item(9, 9).Value = 1.629720535991731
$ws.Cells.Item(9, 10).Value = 0.5200037917938179
$ws.Cells.Item(9, 11).Value = 1.074862163892774
$ws.Cells.Item(9, 12).Value = 0.6909648966888556
$ws.Cells.Item(9, 13).Value = 0.6325680684516349
$ws.Cells.Item(9, 14).Value = 1.031914029066728
$ws.Cells.Item(9, 15).Value = 0.6594978007220411
$ws.Cells.Item(9, 16).Value = 139.8318697842621
$ws.Cells.Item(9, 17).Value = 223.9343017001679

# Row 10
$ws.Cells.Item(10, 1).Value = "model_12_5_8"
$ws.Cells.Item(10, 2).Value = 0.9373156232391298
$ws.Cells.Item(10, 3).Value = 0.6730041514259022
$ws.Cells.Item(10, 4).Value = 0.3814875901847792
$ws.Cells.Item(10, 5).Value = 0.9585342894956631
$ws.Cells.Item(10, 6).Value = 0.8545161615168169
$ws.Cells.Item(10, 7).Value = 0.4191707160349539
$ws.Cells.Item(10, 8).Value = 2.186622745092428
$ws.Cells.Item(10, 9).Value = 1.818032148830124
$ws.Cells.Item(10, 10).Value = 0.5180279040078586
$ws.Cells.Item(10, 11).Value = 1.168030026418991
$ws.Cells.Item(10, 12).Value = 0.6600553509580308
$ws.Cells.Item(10, 13).Value = 0.6474339472370552
$ws.Cells.Item(10, 14).Value = 1.033431667605798
$ws.Cells.Item(10, 15).Value = 0.6749965507439681
$ws.Cells.Item(10, 16).Value = 139.7389540104045
$ws.Cells.Item(10, 17).Value = 223.8413859263104

# Row 11
$ws.Cells.Item(11, 1).Value = "model_12_5_9"
$ws.Cells.Item(11, 2).Value = 0.9340550445674747
$ws.Cells.Item(11, 3).Value = 0.6685499983507529
$ws.Cells.Item(11, 4).Value = 0.3181980075705199
$ws.Cells.Item(11, 5).Value = 0.9581622439736354
$ws.Cells.Item(11, 6).Value = 0.8426411642500742
$ws.Cells.Item(11, 7).Value = 0.4409742206258961
$ws.Cells.Item(11, 8).Value = 2.21640768721544
$ws.Cells.Item(11, 9).Value = 2.004063171090677
$ws.Cells.Item(11, 10).Value = 0.5226758398475529
$ws.Cells.Item(11, 11).Value = 1.263369505469115
$ws.Cells.Item(11, 12).Value = 0.6318767928978746
$ws.Cells.Item(11, 13).Value = 0.6640588984614965
$ws.Cells.Item(11, 14).Value = 1.035170642897347
$ws.Cells.Item(11, 15).Value = 0.6923292605604272
$ws.Cells.Item(11, 16).Value = 139.6375377237491
$ws.Cells.Item(11, 17).Value = 223.739969639655

# Row 12
$ws.Cells.Item(12, 1).Value = "model_12_5_10"
$ws.Cells.Item(12, 2).Value = 0.9304728481841289
$ws.Cells.Item(12, 3).Value = 0.6638154730450276
$ws.Cells.Item(12, 4).Value = 0.2565529783882562
$ws.Cells.Item(12, 5).Value = 0.9571659447247988
$ws.Cells.Item(12, 6).Value = 0.8305815259105677
$ws.Cells.Item(12, 7).Value = 0.4649283843358273
$ws.Cells.Item(12, 8).Value = 2.248067479735309
$ws.Cells.Item(12, 9).Value = 2.185260254755349
$ws.Cells.Item(12, 10).Value = 0.5351225290604502
$ws.Cells.Item(12, 11).Value = 1.360191391907899
$ws.Cells.Item(12, 12).Value = 0.606465460362917
$ws.Cells.Item(12, 13).Value = 0.6818565716745915
$ws.Cells.Item(12, 14).Value = 1.037081147635131
$ws.Cells.Item(12, 15).Value = 0.7108846175684662
$ws.Cells.Item(12, 16).Value = 139.5317437948738
$ws.Cells.Item(12, 17).Value = 223.6341757107796

# Row 13
$ws.Cells.Item(13, 1).Value = "model_12_5_11"
$ws.Cells.Item(13, 2).Value = 0.9266805904912759
$ws.Cells.Item(13, 3).Value = 0.658946846639081
$ws.Cells.Item(13, 4).Value = 0.1974862199196076
$ws.Cells.Item(13, 5).Value = 0.9555367878310037
$ws.Cells.Item(13, 6).Value = 0.8185014749353468
$ws.Cells.Item(13, 7).Value = 0.4902872289896759
$ws.Cells.Item(13, 8).Value = 2.280624006929838
$ws.Cells.Item(13, 9).Value = 2.358878866312825
$ws.Cells.Item(13, 10).Value = 0.5554754597284142
$ws.Cells.Item(13, 11).Value = 1.45717716302062
$ws.Cells.Item(13, 12).Value = 0.583916388234186
$ws.Cells.Item(13, 13).Value = 0.7002051335070858
$ws.Cells.Item(13, 14).Value = 1.03910368507132
$ws.Cells.Item(13, 15).Value = 0.7300143156649291
$ws.Cells.Item(13, 16).Value = 139.4255277560064
$ws.Cells.Item(13, 17).Value = 223.5279596719123

# Remove now-unused trailing rows 14 and 15 (model_12_5_12, model_12_5_13 dropped)
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Delete()
